$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 13.5333
$ws.Range("E6").Value = 12.3367
$ws.Range("E7").Value = 11.7844
$ws.Range("E8").Value = 13.6716
$ws.Range("E16").Value = 12.23900000000001
$ws.Range("E20").Value = 13.39249999999999
$ws.Range("E21").Value = 12.74349999999999
